$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 409
$ws.Range("E11").Value = 277

$ws.Range("E40").Value = 223
$ws.Range("F40").Value = 100
$ws.Range("H40").Value = 100

$ws.Range("E41").Value = 327
$ws.Range("E42").Value = 297
